$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

$ws.Range("A$row").Value = "'2023-06-05"
$ws.Range("B$row").Value = "'18:32:19"
$ws.Range("C$row").Value = "'Monday"
$ws.Range("D$row").Value = "'23"

$textRange = "A" + $row + ":D" + $row
$ws.Range($textRange).ClearFormats()

$ws.Range("E$row").Value = 121216
$ws.Range("F$row").Value = 134221
$ws.Range("G$row").Value = 159827
$ws.Range("H$row").Value = 130175
$ws.Range("I$row").Value = 174987
$ws.Range("J$row").Value = 112617
$ws.Range("K$row").Value = 200166
$ws.Range("L$row").Value = 219583
$ws.Range("M$row").Value = 172298
$ws.Range("N$row").Value = 119399
$ws.Range("O$row").Value = 38283
$ws.Range("P$row").Value = 34681
$ws.Range("Q$row").Value = 50327
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36658
$ws.Range("T$row").Value = -1
